$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BB (54) mirrors column BA (53) but shifted by one forecast period.
# Row 1 holds the new date header; rows 3-18 repeat the BA value; rows 19-21
# hold newly computed values as per the source diff.

# Header date (row 1)
$ws.Range("BA1").Copy($ws.Range("BB1")) | Out-Null
$ws.Range("BB1").Value = 45986

# Rows 3-18: BB value equals the existing BA value for that row
$copyRows = 3..18
foreach ($r in $copyRows) {
    $baCell = "BA" + $r
    $bbCell = "BB" + $r
    $ws.Range($bbCell).Value = $ws.Range($baCell).Value2
}

# Rows 19-21: new forecast values that differ from BA column
$ws.Range("BB19").Value = -2.451276118722334
$ws.Range("BB20").Value = -0.8888225292121632
$ws.Range("BB21").Value = -1.723692879931693
